$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Employee data to write starting at row 2 (row 1 is the header: Empleado - Identificador / Mail)
$data = @(
    ,@("ISMAEL_ASENJO", "ismael.asenjo@stihl.es")
    ,@("ALTAMIRA", "inaki.altamira@stihl.es")
    ,@("ORIOL_PINA", "oriol.pina@stihl.es")
    ,@("TERESA_NIETO", "teresa.nieto@stihl.es")
    ,@("BEGONA_RUIZ", "begona.ruiz@stihl.es")
    ,@("ADOLFOCASASEMPERE", "adolfo.casasempereblanquer@stihl.es")
    ,@("BERND_HULLERUM", "bernd.hullerum@stihl.es")
    ,@("M_MARTINEZ", "m.martinez@stihl.es")
    ,@("MARIO_BALTANAS", "mario.baltanas@stihl.es")
    ,@("PEDRO_MOYANO", "pedro.moyano@stihl.es")
    ,@("FERNANDO_BARRIO", "fernando.barrio@stihl.es")
    ,@("BENJAMIN_MACKH", "benjamin.mackh@stihl.es")
    ,@("ALVARO_LOPEZ", "alvaro.lopez@stihl.es")
    ,@("ANA_SANCHEZ", "ana.sanchez@stihl.es")
    ,@("ENRIQUE_ALONSO", "enrique.alonso@stihl.es")
    ,@("ARANCHA_GANAN", "arancha.ganan@stihl.es")
    ,@("JORGE_HEVIA", "jorge.hevia@stihl.es")
    ,@("VICTOR_PRADA", "victor.prada@stihl.es")
    ,@("RAQUEL_MUNOZ", "raquel.munoz@stihl.es")
    ,@("JSANCHEZ", "jsanchez@stihl.es")
    ,@("JOAQUIN_PORTILLO", "joaquin.portillo@stihl.es")
    ,@("FERNANDO_ALVARO", "fernando.alvarocastillo@stihl.es")
    ,@("ANDRES_LOPEZ", "andres.lopez@stihl.es")
    ,@("JAVIER_SICILIA", "javier.sicilia@stihl.es")
    ,@("TONY_BLANCO", "tony.blanco@stihl.es")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Update selection to reflect the saved view state
$ws.Range("C13").Select()
